$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add a new row of data (row 5) to the BIIB bag sheet, mirroring existing rows.
$ws.Cells.Item(5, 1).Value = 42606.88175925926

$ws.Cells.Item(5, 2).Value = 32
$ws.Cells.Item(5, 3).Value = 68
$ws.Cells.Item(5, 4).Value = 29
$ws.Cells.Item(5, 5).Value = 58
$ws.Cells.Item(5, 6).Value = 41
$ws.Cells.Item(5, 7).Value = 7959
$ws.Cells.Item(5, 8).Value = 13770
$ws.Cells.Item(5, 9).Value = 1678
$ws.Cells.Item(5, 10).Value = 222
$ws.Cells.Item(5, 11).Value = 95
$ws.Cells.Item(5, 12).Value = 7
$ws.Cells.Item(5, 13).Value = 5
$ws.Cells.Item(5, 14).Value = "Bag"

$wb.Save()
